$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header/text cells (A1:C1 keep their text style; D1 is new) ---
# These hold numeric-looking strings ("2", "3", "0", "0.1") that must stay
# stored as TEXT, not be auto-coerced to numbers. We stage each value in a
# scratch cell formatted as Text, then copy only the VALUE (not the format)
# into place so the destination's existing style/format is left untouched
# (A1:C1 keep style index 1; D1 picks it up via an explicit format copy).
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

Set-TextValue "A1" "2"
Set-TextValue "B1" "3"
Set-TextValue "C1" "0"

# D1 is a brand new cell, so give it the same style as the rest of the
# header row (copied from C1) before writing its text value.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
Set-TextValue "D1" "0.1"

# --- Row 2: numeric values (existing row, updated) ---
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

# --- Row 3: numeric values (new row) ---
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
